$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 83
$ws.Range("J33").Value = 49.5
$ws.Range("L33").Value = 49.5
$ws.Range("N33").Value = -507.5
$ws.Range("H74").Value = 31261418
$ws.Range("I74").Value = 62505396
$ws.Range("K74").Value = 62505396
$ws.Range("M74").Value = -62504460
$ws.Range("H77").Value = 31261418
$ws.Range("I77").Value = 62505396
$ws.Range("K77").Value = 312526980
$ws.Range("M77").Value = -312522300
$ws.Range("H132").Value = 581.14813
$ws.Range("I132").Value = 584.96075
$ws.Range("K132").Value = 1754.88225
$ws.Range("M132").Value = 775.1177500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4611
$ws.Range("I45").Value = 3242.0715
$ws.Range("K45").Value = 3242.0715
$ws.Range("M45").Value = -2865.0715
$ws.Range("H61").Value = 52633988
$ws.Range("I61").Value = 1518.8667
$ws.Range("K61").Value = 1518.8667
$ws.Range("M61").Value = -1306.8667
$ws.Range("H74").Value = 30515.666
$ws.Range("I74").Value = 43382.082
$ws.Range("K74").Value = 43382.082
$ws.Range("M74").Value = -42508.082
$ws.Range("H77").Value = 30515.666
$ws.Range("I77").Value = 43382.082
$ws.Range("K77").Value = 216910.41
$ws.Range("M77").Value = -212542.41
$ws.Range("H94").Value = 46710.5
$ws.Range("J94").Value = 46710.5
$ws.Range("L94").Value = 46710.5
$ws.Range("N94").Value = -48512.5
$ws.Range("H98").Value = 38518.5
$ws.Range("J98").Value = 38518.5
$ws.Range("L98").Value = 38518.5
$ws.Range("N98").Value = -44508.5
$ws.Range("H122").Value = 3058.7317
$ws.Range("I122").Value = 2717.4546
$ws.Range("K122").Value = 8152.3638
$ws.Range("M122").Value = -5702.3638
$ws.Range("H136").Value = 52633988
$ws.Range("I136").Value = 1518.8667
$ws.Range("K136").Value = 4556.6001
$ws.Range("M136").Value = -2006.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 75000
$ws.Range("J74").Value = 75000
$ws.Range("L74").Value = 75000
$ws.Range("N74").Value = -76872
$ws.Range("H77").Value = 75000
$ws.Range("J77").Value = 75000
$ws.Range("L77").Value = 225000
$ws.Range("N77").Value = -234360
$ws.Range("H86").Value = 93928.73
$ws.Range("I86").Value = 169204.17
$ws.Range("K86").Value = 169204.17
$ws.Range("M86").Value = -168081.17
$ws.Range("H89").Value = 93928.73
$ws.Range("I89").Value = 169204.17
$ws.Range("K89").Value = 846020.8500000001
$ws.Range("M89").Value = -840404.8500000001
$ws.Range("H134").Value = 5439313.5
$ws.Range("I134").Value = 7815603
$ws.Range("K134").Value = 23446809
$ws.Range("M134").Value = -23444274
$ws.Range("H135").Value = 87936.86
$ws.Range("J135").Value = 87936.86
$ws.Range("L135").Value = 87936.86
$ws.Range("N135").Value = -98076.86
$ws.Range("H139").Value = 61499.668
$ws.Range("J139").Value = 67799.60000000001
$ws.Range("L139").Value = 67799.60000000001
$ws.Range("N139").Value = -78079.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 33930.8
$ws.Range("J28").Value = 33930.8
$ws.Range("L28").Value = 33930.8
$ws.Range("N28").Value = -34420.8
$ws.Range("H31").Value = 5343.952
$ws.Range("I31").Value = 1969.6
$ws.Range("J31").Value = 7218.593
$ws.Range("K31").Value = 1969.6
$ws.Range("L31").Value = 7218.593
$ws.Range("M31").Value = -1674.6
$ws.Range("N31").Value = -7808.593
$ws.Range("H34").Value = 5343.952
$ws.Range("I34").Value = 1969.6
$ws.Range("J34").Value = 7218.593
$ws.Range("K34").Value = 1969.6
$ws.Range("L34").Value = 7218.593
$ws.Range("M34").Value = -1767.6
$ws.Range("N34").Value = -7622.593
$ws.Range("H95").Value = 59164.832
$ws.Range("J95").Value = 59164.832
$ws.Range("L95").Value = 59164.832
$ws.Range("N95").Value = -64656.832
$ws.Range("H122").Value = 1781.625
$ws.Range("I122").Value = 1253.2307
$ws.Range("K122").Value = 3759.6921
$ws.Range("M122").Value = -1309.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1826.5714
$ws.Range("I98").Value = 503
$ws.Range("J98").Value = 2047.1666
$ws.Range("K98").Value = 1509
$ws.Range("L98").Value = 6141.4998
$ws.Range("M98").Value = -11
$ws.Range("N98").Value = -9137.4998
$ws.Range("H122").Value = 3536529.8
$ws.Range("J122").Value = 1053.3334
$ws.Range("L122").Value = 9480.000599999999
$ws.Range("N122").Value = -14380.0006
$ws.Range("H131").Value = 40076.81
$ws.Range("I131").Value = 1419.5
$ws.Range("J131").Value = 43298.25
$ws.Range("K131").Value = 4258.5
$ws.Range("L131").Value = 129894.75
$ws.Range("M131").Value = 781.5
$ws.Range("N131").Value = -139974.75
$ws.Range("H137").Value = 108263.48
$ws.Range("I137").Value = 112636.22
$ws.Range("J137").Value = 104328
$ws.Range("K137").Value = 337908.66
$ws.Range("L137").Value = 312984
$ws.Range("M137").Value = -332808.66
$ws.Range("N137").Value = -323184
$ws.Range("H140").Value = 66359.484
$ws.Range("I140").Value = 101257.2
$ws.Range("K140").Value = 303771.6
$ws.Range("M140").Value = -298591.6
$ws.Range("H141").Value = 6461.5264
$ws.Range("I141").Value = 4409.9165
$ws.Range("J141").Value = 9978.571
$ws.Range("K141").Value = 13229.7495
$ws.Range("L141").Value = 29935.713
$ws.Range("M141").Value = -8049.749500000002
$ws.Range("N141").Value = -40295.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 2499.6667
$ws.Range("I29").Value = 999
$ws.Range("J29").Value = 3250
$ws.Range("K29").Value = 999
$ws.Range("L29").Value = 3250
$ws.Range("M29").Value = -709
$ws.Range("N29").Value = -3830
$ws.Range("H100").Value = 48484.5
$ws.Range("J100").Value = 48484.5
$ws.Range("L100").Value = 48484.5
$ws.Range("N100").Value = -50648.5
$ws.Range("H102").Value = 6292
$ws.Range("I102").Value = 6181
$ws.Range("K102").Value = 6181
$ws.Range("M102").Value = -4559

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 4561.625
$ws.Range("I107").Value = 4561.625
$ws.Range("K107").Value = 4561.625
$ws.Range("M107").Value = -2641.625
$ws.Range("H132").Value = 12504710
$ws.Range("I132").Value = 26318232
$ws.Range("K132").Value = 78954696
$ws.Range("M132").Value = -78952166
$ws.Range("H140").Value = 76416.164
$ws.Range("J140").Value = 76416.164
$ws.Range("L140").Value = 76416.164
$ws.Range("N140").Value = -86776.164

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 337548.66
$ws.Range("I122").Value = 502325.5
$ws.Range("J122").Value = 7995
$ws.Range("K122").Value = 1506976.5
$ws.Range("L122").Value = 23985
$ws.Range("M122").Value = -1504526.5
$ws.Range("N122").Value = -28885
